# Apply the mapping-corps update for fr-lm-technique-imagerie:
#  - bump the "Date" metadata value
#  - rename the "topographique" sub-element to "precisionTopographique"
#    throughout the Elements sheet (columns A, B and AF)
#  - adjust the best-fit column widths that result from the longer text

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-02-04T10:58:36+00:00"

# --- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 13: fr-lm-technique-imagerie.lateralite.topographique -> ...precisionTopographique
$elements.Range("A13").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique"
$elements.Range("B13").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique"
$elements.Range("AF13").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique"

# Row 14: ...topographique.id -> ...precisionTopographique.id
$elements.Range("A14").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.id"
$elements.Range("B14").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.id"

# Row 15: ...topographique.extension -> ...precisionTopographique.extension
$elements.Range("A15").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.extension"
$elements.Range("B15").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.extension"

# Row 16: ...topographique.coding -> ...precisionTopographique.coding
$elements.Range("A16").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.coding"
$elements.Range("B16").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.coding"

# Row 17: ...topographique.text -> ...precisionTopographique.text
$elements.Range("A17").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.text"
$elements.Range("B17").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.text"

# Row 18: ...topographique.nom -> ...precisionTopographique.nom
$elements.Range("A18").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.nom"
$elements.Range("B18").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.nom"
$elements.Range("AF18").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.nom"

# Row 19: ...topographique.valeur -> ...precisionTopographique.valeur
$elements.Range("A19").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.valeur"
$elements.Range("B19").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.valeur"
$elements.Range("AF19").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.valeur"

# --- Column widths (best-fit grew because of the longer strings) -----
# (target stored widths: 55.7890625 for columns A/B, 53.0390625 for AF;
#  the host's ColumnWidth setter quantises to its own pixel grid, so these
#  inputs land on the closest attainable grid value to the target width)
$elements.Columns.Item(1).ColumnWidth = 55
$elements.Columns.Item(2).ColumnWidth = 55
$elements.Columns.Item(32).ColumnWidth = 52.15
